# Update the AC (Final_Scores) column values for existing rows 2-113
$acUpdates = @{
    2 = 4.210556481126316
    3 = 3.909895722116946
    4 = 3.695337173650969
    5 = 3.838383078751985
    6 = 6.828378372877175
    7 = 6.086697918047746
    8 = 6.232130758331659
    9 = 7.203942109422554
    10 = 6.708039749468196
    11 = 6.169641932022943
    12 = 3.722322220335927
    13 = 3.673231583640857
    14 = 6.834595556599855
    15 = 4.848975541471533
    16 = 4.450346307838432
    17 = 4.875435759982659
    18 = 3.676116834841885
    19 = 4.136121377979959
    20 = 6.706682359940232
    21 = 6.416861987488036
    22 = 6.983575932577238
    23 = 6.428886148937228
    24 = 3.784265585097816
    25 = 7.135734103078507
    26 = 5.80769364641284
    27 = 3.934666470541204
    28 = 5.22595302578463
    29 = 5.703339386727711
    30 = 5.191366243965748
    31 = 4.788363843523755
    32 = 4.946177306283944
    33 = 5.857996245439912
    34 = 3.965614561003554
    35 = 4.83575302325406
    36 = 5.977973324415419
    37 = 6.304302388920509
    38 = 7.439842604892244
    39 = 5.924148177686354
    40 = 6.074985769619514
    41 = 5.780083530488447
    42 = 5.267078630486325
    43 = 4.636608631874995
    44 = 5.637078138950337
    45 = 4.423665228815809
    46 = 5.347490534131101
    47 = 6.110264401016741
    48 = 5.882931701281661
    49 = 5.56948303250207
    50 = 6.055059996397329
    51 = 5.857668916869738
    52 = 5.682759147642267
    53 = 5.159511682181955
    54 = 5.823882639033572
    55 = 6.067805355323006
    56 = 5.681858920048636
    57 = 5.748773984490281
    58 = 5.868345406239691
    59 = 5.875649276275656
    60 = 5.772370510241812
    61 = 5.684048590699295
    62 = 6.88114723571185
    63 = 6.726563284170673
    64 = 5.929262969555773
    65 = 6.089455796111457
    66 = 6.518858170534222
    67 = 5.089561283258011
    68 = 5.336657077195814
    69 = 5.558307096702853
    70 = 5.766109690518703
    71 = 4.807913876969971
    72 = 7.987934718712288
    73 = 5.114018997689896
    74 = 5.164282718870368
    75 = 5.203642832574967
    76 = 5.956988646227575
    77 = 5.411660473860315
    78 = 7.781635816195867
    79 = 4.966999288577204
    80 = 4.76963244787856
    81 = 5.860564717361457
    82 = 6.673934739257777
    83 = 4.541809156308467
    84 = 5.802597336979916
    85 = 7.138085861762434
    86 = 5.916056568413392
    87 = 6.61612822795899
    88 = 4.965896600313615
    89 = 4.423561601949645
    90 = 4.450491917916222
    91 = 3.424828500481773
    92 = 6.208741055422141
    93 = 5.647850531824542
    94 = 5.84792104688896
    95 = 5.009639929046132
    96 = 5.687074126322987
    97 = 4.953191580972033
    98 = 5.625920337168784
    99 = 5.896963223304297
    100 = 5.458457410138247
    101 = 4.976008653112045
    102 = 5.630688314236976
    103 = 7.539769353797888
    104 = 6.112787343966355
    105 = 7.075846110082174
    106 = 5.462148723919706
    107 = 3.961228699695829
    108 = 6.099282282992595
    109 = 5.968479705860054
    110 = 4.713129710404051
    111 = 6.111126864155356
    112 = 6.048302447522063
    113 = 5.774504519512302
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $acUpdates.Keys) {
    $ws.Cells.Item($row, 29).Value = $acUpdates[$row]
}

# Column O (Washington DC score) fix for row 10
$ws.Cells.Item(10, 15).Value = 10

# Column N (San Francisco score) fix for row 72
$ws.Cells.Item(72, 14).Value = 10

# Rename the two city entries that were mislabeled with the wrong state
$ws.Cells.Item(112, 3).Value = "Milwaukee – Wisconsin"
$ws.Cells.Item(113, 3).Value = "Madison – Wisconsin"

# Append the new Washington D.C. self-city row (row 114)
$ws.Cells.Item(113, 1).Copy() | Out-Null
$ws.Cells.Item(114, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(114, 1).Value = 112
$ws.Cells.Item(114, 2).Value = "Washington D.C. "
$ws.Cells.Item(114, 3).Value = "Washington D.C."
$ws.Cells.Item(114, 4).Value = $true
$ws.Cells.Item(114, 5).Value = 25
$ws.Cells.Item(114, 6).Value = 20
$ws.Cells.Item(114, 7).Value = 74
$ws.Cells.Item(114, 8).Value = 77
$ws.Cells.Item(114, 9).Value = 28.25
$ws.Cells.Item(114, 10).Value = 103
$ws.Cells.Item(114, 11).Value = 745
$ws.Cells.Item(114, 12).Value = 168.7
$ws.Cells.Item(114, 13).Value = 1
$ws.Cells.Item(114, 14).Value = 6
$ws.Cells.Item(114, 15).Value = 2
$ws.Cells.Item(114, 16).Value = 10
$ws.Cells.Item(114, 17).Value = 17
$ws.Cells.Item(114, 18).Value = 63
$ws.Cells.Item(114, 19).Value = 20
$ws.Cells.Item(114, 20).Value = 9
$ws.Cells.Item(114, 21).Value = 68
$ws.Cells.Item(114, 22).Value = 24123
$ws.Cells.Item(114, 23).Value = 642
$ws.Cells.Item(114, 24).Value = 50
$ws.Cells.Item(114, 25).Value = 144
$ws.Cells.Item(114, 26).Value = 899
$ws.Cells.Item(114, 27).Value = 5.75
$ws.Cells.Item(114, 28).Value = 5.75
$ws.Cells.Item(114, 29).Value = 5.930666430938779

Write-Output "edit complete"
